$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "AUTOR" header (K1) becomes "EJECUTIVO" ...
$ws.Range("K1").Value = "EJECUTIVO"

# ... and a brand-new "CLASIFICADOR" header is appended in L1, copying the
# same header formatting (bold red font on yellow fill, centered).
$ws.Range("L1").Value = "CLASIFICADOR"
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null

# Re-apply the AutoFilter so it spans the new column (toggling it off first,
# since re-running AutoFilter on an active filter would just disable it).
$ws.AutoFilterMode = $false
$ws.Range("A1:L1").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name needs to track the new range too.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$L`$1"
    }
}

# Widen K and L to fit their new header text.
$ws.Columns.Item(11).ColumnWidth = 14.998697916666666
$ws.Columns.Item(12).ColumnWidth = 18.166666666666668

# Move the active selection, matching the edit's recorded cursor position.
$ws.Range("L4").Select()
